$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename a handful of cell_name labels in column B to their updated /
# corrected English spellings. Order matters for how new shared strings
# get appended, so apply them row-by-row in the same order as the source
# edit (suburbs/Bnei Brak/Tel Aviv first, then Jerusalem/Jezreel, then
# Beer Sheva Arabs).
$renames = [ordered]@{
    16 = "TLV suburbs"
    17 = "Bnei Brak"
    15 = "Tel Aviv - Yafo"
    2  = "Jerusalem and sub."
    6  = "Jezreel Valley"
    20 = "Beer Sheva Arabs"
}

foreach ($r in $renames.Keys) {
    $ws.Cells.Item($r, 2).Value = $renames[$r]
}

# Update the active selection recorded in the sheet view.
$ws.Range("D14").Select()
